# Auto-generated Excel COM-interop edits for Lich_Profits workbook
# Updates market-price derived values (H..N columns) across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 188.14285
$ws.Range("I2").Value = 207.75
$ws.Range("J2").Value = 162
$ws.Range("K2").Value = 207.75
$ws.Range("L2").Value = 162
$ws.Range("M2").Value = -94.75
$ws.Range("N2").Value = -388

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 14432.5
$ws.Range("I51").Value = 17664.834
$ws.Range("J51").Value = 13550.954
$ws.Range("K51").Value = 17664.834
$ws.Range("L51").Value = 13550.954
$ws.Range("M51").Value = -17180.834
$ws.Range("N51").Value = -14518.954

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 50002628
$ws.Range("I62").Value = 33334438
$ws.Range("K62").Value = 33334438
$ws.Range("M62").Value = -33333814

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 50002628
$ws.Range("I65").Value = 33334438
$ws.Range("K65").Value = 166672190
$ws.Range("M65").Value = -166669070

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 4009.8
$ws.Range("J88").Value = 4899.75
$ws.Range("L88").Value = 4899.75
$ws.Range("N88").Value = -5711.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 4009.8
$ws.Range("J91").Value = 4899.75
$ws.Range("L91").Value = 4899.75
$ws.Range("N91").Value = -7707.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 5363319
$ws.Range("J112").Value = 6338084.5
$ws.Range("L112").Value = 19014253.5
$ws.Range("N112").Value = -19016469.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 333335330
$ws.Range("I127").Value = 200002400
$ws.Range("K127").Value = 600007200
$ws.Range("M127").Value = -600002240

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 350
$ws.Range("I129").Value = 350
$ws.Range("K129").Value = 1050
$ws.Range("M129").Value = 3950

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4047.6924
$ws.Range("I132").Value = 3863.6667
$ws.Range("K132").Value = 11591.0001
$ws.Range("M132").Value = -9061.000100000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2386.4607
$ws.Range("J138").Value = 3386.6345
$ws.Range("L138").Value = 10159.9035
$ws.Range("N138").Value = -20439.9035

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2578
$ws.Range("I141").Value = 2587.3462
$ws.Range("J141").Value = 2517.25
$ws.Range("K141").Value = 7762.0386
$ws.Range("L141").Value = 7551.75
$ws.Range("M141").Value = -2582.0386
$ws.Range("N141").Value = -17911.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 16005
$ws.Range("I39").Value = 3015
$ws.Range("J39").Value = 22500
$ws.Range("K39").Value = 3015
$ws.Range("L39").Value = 22500
$ws.Range("M39").Value = -2495
$ws.Range("N39").Value = -23540

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 90895.09
$ws.Range("I74").Value = 94842.86
$ws.Range("J74").Value = 7992
$ws.Range("K74").Value = 94842.86
$ws.Range("L74").Value = 7992
$ws.Range("M74").Value = -93968.86
$ws.Range("N74").Value = -9740

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 90895.09
$ws.Range("I77").Value = 94842.86
$ws.Range("J77").Value = 7992
$ws.Range("K77").Value = 474214.3
$ws.Range("L77").Value = 39960
$ws.Range("M77").Value = -469846.3
$ws.Range("N77").Value = -48696

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3213.3333
$ws.Range("J88").Value = 3830
$ws.Range("L88").Value = 3830
$ws.Range("N88").Value = -4642

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 3213.3333
$ws.Range("J91").Value = 3830
$ws.Range("L91").Value = 3830
$ws.Range("N91").Value = -6638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2470.5334
$ws.Range("I102").Value = 2389.6924
$ws.Range("K102").Value = 2389.6924
$ws.Range("M102").Value = -767.6923999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1808.0328
$ws.Range("I132").Value = 1637.0465
$ws.Range("K132").Value = 4911.139499999999
$ws.Range("M132").Value = -2381.139499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 49250
$ws.Range("J139").Value = 49250
$ws.Range("L139").Value = 49250
$ws.Range("N139").Value = -59530

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1413.12
$ws.Range("I20").Value = 1264.091
$ws.Range("K20").Value = 1264.091
$ws.Range("M20").Value = -1017.091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 287402.16
$ws.Range("I31").Value = 386181.06
$ws.Range("K31").Value = 386181.06
$ws.Range("M31").Value = -385886.06

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 287402.16
$ws.Range("I34").Value = 386181.06
$ws.Range("K34").Value = 386181.06
$ws.Range("M34").Value = -385979.06

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2716.389
$ws.Range("I58").Value = 2452.2666
$ws.Range("K58").Value = 2452.2666
$ws.Range("M58").Value = -2249.2666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 36143.08
$ws.Range("J68").Value = 36143.08
$ws.Range("L68").Value = 36143.08
$ws.Range("N68").Value = -37641.08

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 36143.08
$ws.Range("J71").Value = 36143.08
$ws.Range("L71").Value = 108429.24
$ws.Range("N71").Value = -115917.24

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3408.2083
$ws.Range("I132").Value = 3461.7856
$ws.Range("J132").Value = 3333.2
$ws.Range("K132").Value = 10385.3568
$ws.Range("L132").Value = 9999.599999999999
$ws.Range("M132").Value = -7855.356800000001
$ws.Range("N132").Value = -15059.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5596.8184
$ws.Range("I134").Value = 6010.185
$ws.Range("J134").Value = 3736.6667
$ws.Range("K134").Value = 18030.555
$ws.Range("L134").Value = 11210.0001
$ws.Range("M134").Value = -15495.555
$ws.Range("N134").Value = -16280.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2716.389
$ws.Range("I136").Value = 2452.2666
$ws.Range("K136").Value = 7356.7998
$ws.Range("M136").Value = -4806.7998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 574.2917
$ws.Range("I39").Value = 349.2857
$ws.Range("K39").Value = 1047.8571
$ws.Range("M39").Value = -753.8571000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 16667687
$ws.Range("J68").Value = 1433
$ws.Range("L68").Value = 4299
$ws.Range("N68").Value = -5921

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 16667687
$ws.Range("J71").Value = 1433
$ws.Range("L71").Value = 12897
$ws.Range("N71").Value = -21009

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 338.2143
$ws.Range("J98").Value = 369.6
$ws.Range("L98").Value = 1108.8
$ws.Range("N98").Value = -4104.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1079.2222
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 1089.125
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 3267.375
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -7607.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 1552.6666
$ws.Range("I123").Value = 1552.6666
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 4657.9998
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -2207.9998
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 134921940
$ws.Range("I129").Value = 138890420
$ws.Range("J129").Value = 111111110
$ws.Range("K129").Value = 416671260
$ws.Range("L129").Value = 333333330
$ws.Range("M129").Value = -416666260
$ws.Range("N129").Value = -333343330

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 15626427
$ws.Range("I131").Value = 125000910
$ws.Range("J131").Value = 1500.4108
$ws.Range("K131").Value = 375002730
$ws.Range("L131").Value = 4501.232400000001
$ws.Range("M131").Value = -374997690
$ws.Range("N131").Value = -14581.2324

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2303.261
$ws.Range("I139").Value = 1490.7646
$ws.Range("K139").Value = 4472.293799999999
$ws.Range("M139").Value = 667.7062000000005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 965.25
$ws.Range("I31").Value = 965.25
$ws.Range("K31").Value = 965.25
$ws.Range("M31").Value = -673.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H37").Value = 965.25
$ws.Range("I37").Value = 965.25
$ws.Range("K37").Value = 965.25
$ws.Range("M37").Value = -688.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5407.3115
$ws.Range("I70").Value = 5271.838
$ws.Range("J70").Value = 5616.1665
$ws.Range("K70").Value = 5271.838
$ws.Range("L70").Value = 5616.1665
$ws.Range("M70").Value = -5001.838
$ws.Range("N70").Value = -6156.1665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5407.3115
$ws.Range("I73").Value = 5271.838
$ws.Range("J73").Value = 5616.1665
$ws.Range("K73").Value = 5271.838
$ws.Range("L73").Value = 5616.1665
$ws.Range("M73").Value = -4335.838
$ws.Range("N73").Value = -7488.1665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5379.143
$ws.Range("I102").Value = 5379.143
$ws.Range("K102").Value = 5379.143
$ws.Range("M102").Value = -3757.143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1117.125
$ws.Range("I22").Value = 856.4375
$ws.Range("J22").Value = 1377.8125
$ws.Range("K22").Value = 856.4375
$ws.Range("L22").Value = 1377.8125
$ws.Range("M22").Value = -561.4375
$ws.Range("N22").Value = -1967.8125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1117.125
$ws.Range("I27").Value = 856.4375
$ws.Range("J27").Value = 1377.8125
$ws.Range("K27").Value = 856.4375
$ws.Range("L27").Value = 1377.8125
$ws.Range("M27").Value = -749.4375
$ws.Range("N27").Value = -1591.8125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 2419.5
$ws.Range("J31").Value = 5258.5
$ws.Range("L31").Value = 5258.5
$ws.Range("N31").Value = -5754.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 3934
$ws.Range("J31").Value = 4651
$ws.Range("L31").Value = 4651
$ws.Range("N31").Value = -5347

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1514.2142
$ws.Range("I132").Value = 1432.8334
$ws.Range("J132").Value = 2002.5
$ws.Range("K132").Value = 4298.5002
$ws.Range("L132").Value = 6007.5
$ws.Range("M132").Value = -1768.5002
$ws.Range("N132").Value = -11067.5
